# This script reorganizes the "EDA" deck so that instead of dumping every
# chart image + the full statistics table onto slide 1, each chart gets its
# own slide (2-15), and slide 1 becomes a "Statistical Analysis" slide that
# uses a Title+Content layout for the stats table.

function Find-PictureByAltText($slide, $altText) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Type -eq 13 -and $shp.AlternativeText -eq $altText) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# --- Step 1: slide 1 title text -------------------------------------------
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Statistical Analysis - univariate_analysis.txt"

# --- Step 2: copy each chart picture from slide 1 onto its own new slide --
# (pictures are located dynamically by their original alt-text/descr so the
# later layout / shape-count changes on slide 1 cannot shift indices out
# from under us)
$chartSlides = @(
    @{ Slide = 2;  Image = "univariate_bmi.png";               Title = "Chart - univariate_bmi.png" },
    @{ Slide = 3;  Image = "univariate_charges.png";            Title = "Chart - univariate_charges.png" },
    @{ Slide = 4;  Image = "univariate_children.png";           Title = "Chart - univariate_children.png" },
    @{ Slide = 5;  Image = "univariate_region.png";              Title = "Chart - univariate_region.png" },
    @{ Slide = 6;  Image = "univariate_sex.png";                Title = "Chart - univariate_sex.png" },
    @{ Slide = 7;  Image = "univariate_smoker.png";             Title = "Chart - univariate_smoker.png" },
    @{ Slide = 8;  Image = "bivariate_age_vs_bmi.png";          Title = "Chart - bivariate_age_vs_bmi.png" },
    @{ Slide = 9;  Image = "bivariate_age_vs_charges.png";      Title = "Chart - bivariate_age_vs_charges.png" },
    @{ Slide = 10; Image = "bivariate_age_vs_children.png";     Title = "Chart - bivariate_age_vs_children.png" },
    @{ Slide = 11; Image = "bivariate_bmi_vs_charges.png";      Title = "Chart - bivariate_bmi_vs_charges.png" },
    @{ Slide = 12; Image = "bivariate_bmi_vs_children.png";     Title = "Chart - bivariate_bmi_vs_children.png" },
    @{ Slide = 13; Image = "bivariate_children_vs_charges.png"; Title = "Chart - bivariate_children_vs_charges.png" },
    @{ Slide = 14; Image = "correlation_matrix.png";            Title = "Chart - correlation_matrix.png" },
    @{ Slide = 15; Image = "univariate_age.png";                Title = "Chart - univariate_age.png" }
)

foreach ($item in $chartSlides) {
    $targetSlide = $p.Slides.Item($item.Slide)
    $targetSlide.Shapes.Item(1).TextFrame.TextRange.Text = $item.Title

    $srcPic = Find-PictureByAltText $s1 $item.Image
    $srcPic.Copy()
    $newPic = $targetSlide.Shapes.Paste().Item(1)
    $newPic.Name = "Picture 2"
    $newPic.Top = 108
    $newPic.Left = 72
}

# --- Step 3: strip all the chart pictures back off of slide 1 -------------
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Type -eq 13) {
        $shp.Delete()
    }
}

# --- Step 4: give slide 1 a Content Placeholder by switching its layout ---
# (the master's 2nd custom layout is "Title and Content", which owns the
# "Content Placeholder 2" shape with <p:ph idx="1"/> that the edit adds)
$sm = $p.SlideMaster
$titleAndContent = $sm.CustomLayouts.Item(2)
$s1.CustomLayout = $titleAndContent

# move the freshly added placeholder up so it sits right after the title
# (send-to-back puts it first, then bring-forward once puts it 2nd)
$contentPh = $s1.Shapes.Item($s1.Shapes.Count)
$contentPh.ZOrder(1)
$contentPh.ZOrder(2)

# --- Step 5: reposition / rename the stats TextBox -------------------------
$statsBox = $s1.Shapes.Item(3)
$statsBox.Name = "TextBox 3"
$statsBox.Top = 108
